$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("returned")

$ws.Range("A3").Value = "АДМ"
$ws.Range("B3").Value = "Адмирал"
$ws.Range("C3").Value = "Дарьин Александр"
$ws.Range("D3").Value = "1369_АДМ_дарьиналександр"
$ws.Range("A4").Value = "АДМ"
$ws.Range("B4").Value = "Адмирал"
$ws.Range("C4").Value = "Педан Руслан"
$ws.Range("D4").Value = "1369_АДМ_педанруслан"
$ws.Range("A5").Value = "АКБ"
$ws.Range("B5").Value = "Ак Барс"
$ws.Range("C5").Value = "Яруллин Альберт"
$ws.Range("D5").Value = "1369_АКБ_яруллинальберт"
$ws.Range("A6").Value = "АМР"
$ws.Range("B6").Value = "Амур"
$ws.Range("C6").Value = "Абросимов Роман"
$ws.Range("D6").Value = "1369_АМР_абросимовроман"
$ws.Range("A7").Value = "БАР"
$ws.Range("B7").Value = "Барыс"
$ws.Range("C7").Value = "Уотерспун Тайлер"
$ws.Range("D7").Value = "1369_БАР_уотерспунтайлер"
$ws.Range("A8").Value = "ДИН"
$ws.Range("B8").Value = "Динамо М"
$ws.Range("C8").Value = "Готовец Кирилл"
$ws.Range("D8").Value = "1369_ДИН_готовецкирилл"
$ws.Range("A9").Value = "ЛАД"
$ws.Range("B9").Value = "Лада"
$ws.Range("C9").Value = "Обидин Андрей"
$ws.Range("D9").Value = "1369_ЛАД_обидинандрей"
$ws.Range("A10").Value = "ЛАД"
$ws.Range("B10").Value = "Лада"
$ws.Range("C10").Value = "Ожгихин Алексей"
$ws.Range("D10").Value = "1369_ЛАД_ожгихиналексей"
$ws.Range("A11").Value = "ЛОК"
$ws.Range("B11").Value = "Локомотив"
$ws.Range("C11").Value = "Сергеев Андрей"
$ws.Range("D11").Value = "1369_ЛОК_сергеевандрей"
$ws.Range("A12").Value = "ММГ"
$ws.Range("B12").Value = "Металлург Мг"
$ws.Range("C12").Value = "Козлов Андрей Е"
$ws.Range("D12").Value = "1369_ММГ_козловандрейе"
$ws.Range("A13").Value = "ММГ"
$ws.Range("B13").Value = "Металлург Мг"
$ws.Range("C13").Value = "Сиряцкий Александр"
$ws.Range("D13").Value = "1369_ММГ_сиряцкийалександр"
$ws.Range("A14").Value = "НХК"
$ws.Range("B14").Value = "Нефтехимик"
$ws.Range("C14").Value = "Хлыстов Никита"
$ws.Range("D14").Value = "1369_НХК_хлыстовникита"
$ws.Range("A15").Value = "СЕВ"
$ws.Range("B15").Value = "Северсталь"
$ws.Range("C15").Value = "Фомин Макар"
$ws.Range("D15").Value = "1369_СЕВ_фоминмакар"
$ws.Range("A17").Value = "СИБ"
$ws.Range("B17").Value = "Сибирь"
$ws.Range("C17").Value = "Аланов Егор"
$ws.Range("D17").Value = "1369_СИБ_алановегор"
$ws.Range("A18").Value = "СКА"
$ws.Range("B18").Value = "СКА"
$ws.Range("C18").Value = "Зайцев Никита И"
$ws.Range("D18").Value = "1369_СКА_зайцевникитаи"
$ws.Range("A19").Value = "СКА"
$ws.Range("B19").Value = "СКА"
$ws.Range("C19").Value = "Короткий Матвей"
$ws.Range("D19").Value = "1369_СКА_короткийматвей"
$ws.Range("A20").Value = "СОЧ"
$ws.Range("B20").Value = "ХК Сочи"
$ws.Range("C20").Value = "Бикмуллин Рафаэль"
$ws.Range("D20").Value = "1369_СОЧ_бикмуллинрафаэль"
$ws.Range("A21").Value = "СОЧ"
$ws.Range("B21").Value = "ХК Сочи"
$ws.Range("C21").Value = "Венгрыжановский Денис"
$ws.Range("D21").Value = "1369_СОЧ_венгрыжановскийденис"
$ws.Range("A22").Value = "СОЧ"
$ws.Range("B22").Value = "ХК Сочи"
$ws.Range("C22").Value = "Хёфенмайер Ноэль"
$ws.Range("D22").Value = "1369_СОЧ_хефенмайерноэль"
$ws.Range("A23").Value = "СПР"
$ws.Range("B23").Value = "Спартак"
$ws.Range("C23").Value = "Вишневский Дмитрий"
$ws.Range("D23").Value = "1369_СПР_вишневскийдмитрий"
$ws.Range("A24").Value = "ТОР"
$ws.Range("B24").Value = "Торпедо"
$ws.Range("C24").Value = "Кручинин Алексей"
$ws.Range("D24").Value = "1369_ТОР_кручининалексей"
$ws.Range("A25").Value = "ТРК"
$ws.Range("B25").Value = "Трактор"
$ws.Range("C25").Value = "Мыльников Сергей И"
$ws.Range("D25").Value = "1369_ТРК_мыльниковсергейи"
$ws.Range("A26").Value = "ТРК"
$ws.Range("B26").Value = "Трактор"
$ws.Range("C26").Value = "Светлаков Андрей"
$ws.Range("D26").Value = "1369_ТРК_светлаковандрей"
$ws.Range("A27").Value = "ЦСК"
$ws.Range("B27").Value = "ЦСКА"
$ws.Range("C27").Value = "Бучельников Дмитрий"
$ws.Range("D27").Value = "1369_ЦСК_бучельниковдмитрий"
$ws.Range("A28").Value = "ЦСК"
$ws.Range("B28").Value = "ЦСКА"
$ws.Range("C28").Value = "Моисеев Данила"
$ws.Range("D28").Value = "1369_ЦСК_моисеевданила"
$ws.Range("A29").Value = "ШДР"
$ws.Range("B29").Value = "Драконы"
$ws.Range("C29").Value = "Бишофф Джейк"
$ws.Range("D29").Value = "1369_ШДР_бишоффджейк"
